# Apply the "LEP" table addition + Natl_Origin reorder edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Reorder the Natl_Origin sheet data rows so SCAG moves from the top
#    (row 2) to the bottom (row 8), matching the alphabetical-then-SCAG
#    ordering used by the other county tables.
# ---------------------------------------------------------------------
$natl = $wb.Worksheets.Item("Natl_Origin")

$natl.Range("A2").Value = "Imperial"
$natl.Range("B2").Value = 29.76
$natl.Range("A3").Value = "Los Angeles"
$natl.Range("B3").Value = 33.51
$natl.Range("A4").Value = "Orange"
$natl.Range("B4").Value = 29.86
$natl.Range("A5").Value = "Riverside"
$natl.Range("B5").Value = 21.55
$natl.Range("A6").Value = "San Bernardino"
$natl.Range("B6").Value = 20.92
$natl.Range("A7").Value = "Ventura"
$natl.Range("B7").Value = 21.2
$natl.Range("A8").Value = "SCAG"
$natl.Range("B8").Value = 29.32

# ---------------------------------------------------------------------
# 2. Add a new "LEP" sheet (after Natl_Origin) with the Limited English
#    Proficiency data table.
# ---------------------------------------------------------------------
$lep = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $natl)
$lep.Name = "LEP"

$lep.Range("A1").Value = "county"
$lep.Range("B1").Value = "lep_perc"
$lep.Range("A1:B1").Font.Bold = $true
$lep.Range("A1:B1").HorizontalAlignment = -4108

$lep.Range("A2").Value = "Imperial"
$lep.Range("B2").Value = 18.31
$lep.Range("A3").Value = "Los Angeles"
$lep.Range("B3").Value = 12.01
$lep.Range("A4").Value = "Orange"
$lep.Range("B4").Value = 8.59
$lep.Range("A5").Value = "Riverside"
$lep.Range("B5").Value = 7.34
$lep.Range("A6").Value = "San Bernardino"
$lep.Range("B6").Value = 7.19
$lep.Range("A7").Value = "Ventura"
$lep.Range("B7").Value = 8.59
$lep.Range("A8").Value = "SCAG"
$lep.Range("B8").Value = 10.18

# ---------------------------------------------------------------------
# 3. Add the new row to the TOC sheet describing the LEP table.
# ---------------------------------------------------------------------
$toc = $wb.Worksheets.Item("TOC")
$toc.Range("A18").Value = "LEP"
$toc.Range("B18").Value = "Limited English Proficiency (%) by County and SCAG Region"

# Restore TOC as the active/selected tab (Worksheets.Add() above made the
# new LEP sheet active, but the original workbook kept TOC selected).
$toc.Activate()
